$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" column (G) for rows 2 and 4 shared the same
# timestamp "2016-08-31 12:17:35" -> "2016-08-31 12:18:24"
$wsOverview.Range("G2").Value = "2016-08-31 12:18:24"
$wsOverview.Range("G4").Value = "2016-08-31 12:18:24"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E) changed from "ht" to "mt" for rows 2 and 4
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime (H) "2016-08-31 12:17:30" -> "2016-08-31 12:18:19"
$wsZhCn.Range("H2").Value = "2016-08-31 12:18:19"
$wsZhCn.Range("H4").Value = "2016-08-31 12:18:19"
# Correspond Handback DateTime (K) "2016-08-31 12:17:47" -> "2016-08-31 12:18:37"
$wsZhCn.Range("K2").Value = "2016-08-31 12:18:37"
$wsZhCn.Range("K4").Value = "2016-08-31 12:18:37"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handback DateTime (K) "2016-08-31 12:17:54" -> "2016-08-31 12:18:44"
$wsDeDe.Range("K2").Value = "2016-08-31 12:18:44"
$wsDeDe.Range("K4").Value = "2016-08-31 12:18:44"
